$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150, pushing existing rows 150-332 down to 151-333
$ws.Rows("150:150").Insert()

# Populate the newly inserted row 150 with a new data record (same constant columns
# as the rest of the dataset, with a new date and volume/price figures)
$ws.Cells.Item(150, 1).Value = 3
$ws.Cells.Item(150, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(150, 3).Value = "Coquimbo"
$ws.Cells.Item(150, 4).Value = 44740
$ws.Cells.Item(150, 5).Value = 5
$ws.Cells.Item(150, 6).Value = 100112039
$ws.Cells.Item(150, 7).Value = "Ciboulette"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 120
$ws.Cells.Item(150, 11).Value = 1500
$ws.Cells.Item(150, 12).Value = 1500
$ws.Cells.Item(150, 13).Value = 1500
$ws.Cells.Item(150, 14).Value = "$/docena de atados"
$ws.Cells.Item(150, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(150, 16).Value = 500
$ws.Cells.Item(150, 17).Value = 3
$ws.Cells.Item(150, 18).Value = "Hortaliza"
